# SCD0016-024 workbook update
# - rename sheet SCD0250 -> SCD0016
# - TC_ID cell (B2) changes from "DGS-265" to "SCD0016-024"
# - column B widened to fit the new, longer TC_ID text
# - selection moves from N1 (scrolled to J1) back to B3 (scrolled to A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "SCD0016"

# Update the TC_ID value in B2
$ws.Range("B2").Value = "SCD0016-024"

# Widen column B to fit the new TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6

# Reset the view: scroll back to A1 and select B3
$ws.Range("A1").Select() | Out-Null
$ws.Range("B3").Select() | Out-Null
